$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text formatting to D:E columns for data rows so that
# numeric-looking price/volume strings are stored as text, matching
# the original workbook (which stores these as inline strings),
# and keep the underlying cell style unchanged (reset to Normal after).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '25.003.25'
$ws.Range("E2").Value = '  +3.02%  '

$ws.Range("D3").Value = '1.729.32'
$ws.Range("E3").Value = '  +3.01%  '

$ws.Range("D4").Value = '1.013'
$ws.Range("E4").Value = '  +1.11%  '

$ws.Range("D5").Value = '314.45'
$ws.Range("E5").Value = '  +2.34%  '

$ws.Range("D6").Value = '1.009'
$ws.Range("E6").Value = '  +1.24%  '

$ws.Range("D7").Value = '0.3801'
$ws.Range("E7").Value = '  +2.09%  '

$ws.Range("D8").Value = '0.3537'
$ws.Range("E8").Value = '  +2.68%  '

$ws.Range("D9").Value = '49.68'
$ws.Range("E9").Value = '  +3.04%  '

$ws.Range("D10").Value = '1.204'
$ws.Range("E10").Value = '  +1.43%  '

$ws.Range("D11").Value = '0.07550'
$ws.Range("E11").Value = '  +3.57%  '

$ws.Range("D12").Value = '1.009'
$ws.Range("E12").Value = '  +1.06%  '

$ws.Range("D13").Value = '6.427'
$ws.Range("E13").Value = '  +5.32%  '

$ws.Range("D14").Value = '21.06'
$ws.Range("E14").Value = '  +2.26%  '

$ws.Range("D15").Value = '7.043'
$ws.Range("E15").Value = '  +4.02%  '

$ws.Range("D16").Value = '1.737.83'
$ws.Range("E16").Value = '  +3.44%  '

$ws.Range("D17").Value = '0.00001136'
$ws.Range("E17").Value = '  +2.09%  '

$ws.Range("D18").Value = '1.010'
$ws.Range("E18").Value = '  +1.28%  '

$ws.Range("D19").Value = '0.06724'
$ws.Range("E19").Value = '  +0.26%  '

$ws.Range("D20").Value = '85.33'
$ws.Range("E20").Value = '  +4.38%  '

$ws.Range("D21").Value = '17.45'
$ws.Range("E21").Value = '  +5.81%  '

$ws.Range("D22").Value = '6.448'
$ws.Range("E22").Value = '  +5.25%  '

$ws.Range("D23").Value = '13.19'
$ws.Range("E23").Value = '  +9.86%  '

$ws.Range("D24").Value = '24.971.42'
$ws.Range("E24").Value = '  +3.05%  '

$ws.Range("D25").Value = '2.450'
$ws.Range("E25").Value = '  +2.09%  '

$ws.Range("D26").Value = '2.833'
$ws.Range("E26").Value = '  +5.98%  '

$ws.Range("D27").Value = '20.63'
$ws.Range("E27").Value = '  +5.71%  '

$ws.Range("D28").Value = '152.26'
$ws.Range("E28").Value = '  +0.61%  '

$ws.Range("D29").Value = '1.925.60'
$ws.Range("E29").Value = '  +3.45%  '

$ws.Range("D30").Value = '132.54'
$ws.Range("E30").Value = '  +4.26%  '

$ws.Range("D31").Value = '1.193'
$ws.Range("E31").Value = '  +20.07%  '

$ws.Range("D32").Value = '6.939'
$ws.Range("E32").Value = '  +8.12%  '

$ws.Range("D33").Value = '4.248'
$ws.Range("E33").Value = '  +5.21%  '

$ws.Range("B34").Value = 'WEMIXTOKEN'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '1.804'
$ws.Range("E34").Value = '  +3.37%  '

$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = '13.85'
$ws.Range("E35").Value = '  +11.37%  '

$ws.Range("D36").Value = '0.08748'
$ws.Range("E36").Value = '  +3.54%  '

$ws.Range("D37").Value = '5.664'
$ws.Range("E37").Value = '  +5.48%  '

$ws.Range("D38").Value = '0.02481'
$ws.Range("E38").Value = '  +5.69%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '9.240'
$ws.Range("E39").Value = '  +3.77%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '0.06630'
$ws.Range("E40").Value = '  +3.22%  '

$ws.Range("D41").Value = '0.2230'
$ws.Range("E41").Value = '  +5.56%  '

$ws.Range("D42").Value = '1.276'
$ws.Range("E42").Value = '  -1.56%  '

$ws.Range("D43").Value = '0.6516'
$ws.Range("E43").Value = '  +5.80%  '

$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = '1.009'
$ws.Range("E44").Value = '  +1.30%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '13.98'
$ws.Range("E45").Value = '  +5.84%  '

$ws.Range("D46").Value = '0.6219'
$ws.Range("E46").Value = '  +4.03%  '

$ws.Range("D47").Value = '3.880'
$ws.Range("E47").Value = '  +2.06%  '

$ws.Range("D48").Value = '2.164'
$ws.Range("E48").Value = '  +7.01%  '

$ws.Range("D49").Value = '130.29'
$ws.Range("E49").Value = '  +2.16%  '

$ws.Range("D50").Value = '0.07303'
$ws.Range("E50").Value = '  +2.27%  '

$ws.Range("D51").Value = '80.32'
$ws.Range("E51").Value = '  +4.94%  '

# Restore the D:E column cells to the default (unstyled) look so the
# saved file does not pick up a stray "@" number format on these cells.
$ws.Range("D2:E51").Style = "Normal"